$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-coerced to a number by Excel
# are forced to remain text (matching the original inline-string cell type)
# by temporarily applying a text number format, then restoring the default style.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "24.723.13"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "1.696.50"
$ws.Range("E3").Value = "  +0.43%  "
Set-TextValue "D4" "0.9994"
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue "D5" "316.80"
$ws.Range("E6").Value = "  +0.31%  "
Set-TextValue "D7" "0.3953"
$ws.Range("E7").Value = "  +1.12%  "
Set-TextValue "D8" "0.4064"
$ws.Range("E8").Value = "  +0.78%  "
Set-TextValue "D9" "1.493"
$ws.Range("E9").Value = "  +2.09%  "
Set-TextValue "D10" "1.002"
$ws.Range("E10").Value = "  +0.36%  "
Set-TextValue "D11" "52.22"
$ws.Range("E11").Value = "  -3.76%  "
Set-TextValue "D12" "0.08896"
$ws.Range("E12").Value = "  +1.99%  "
Set-TextValue "D13" "7.217"
$ws.Range("E13").Value = "  +0.29%  "
Set-TextValue "D14" "23.62"
$ws.Range("E14").Value = "  +2.89%  "
Set-TextValue "D15" "8.142"
$ws.Range("E15").Value = "  +9.53%  "
Set-TextValue "D16" "0.00001329"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "1.694.28"
$ws.Range("E17").Value = "  +0.40%  "
Set-TextValue "D18" "100.06"
$ws.Range("E18").Value = "  -0.24%  "
Set-TextValue "D19" "0.07009"
$ws.Range("E19").Value = "  +0.13%  "
Set-TextValue "D20" "19.70"
$ws.Range("E20").Value = "  +2.29%  "
Set-TextValue "D21" "7.035"
$ws.Range("E21").Value = "  +5.23%  "
$ws.Range("E22").Value = "  +0.44%  "
Set-TextValue "D23" "14.44"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").Value = "24.664.36"
$ws.Range("E24").Value = "  +1.63%  "
Set-TextValue "D25" "3.227"
$ws.Range("E25").Value = "  +7.71%  "
Set-TextValue "D26" "2.347"
$ws.Range("E26").Value = "  +1.49%  "
Set-TextValue "D27" "22.78"
$ws.Range("E27").Value = "  +2.98%  "
Set-TextValue "D28" "163.40"
$ws.Range("E28").Value = "  +1.78%  "
Set-TextValue "D29" "136.61"
$ws.Range("E29").Value = "  +3.54%  "
Set-TextValue "D30" "5.165"
$ws.Range("E30").Value = "  +1.29%  "
Set-TextValue "D31" "7.470"
$ws.Range("E31").Value = "  -3.45%  "
$ws.Range("D32").Value = "1.881.63"
$ws.Range("E32").Value = "  +0.35%  "
Set-TextValue "D33" "1.070"
$ws.Range("E33").Value = "  -1.22%  "
Set-TextValue "D34" "0.08618"
$ws.Range("E34").Value = "  -0.82%  "
Set-TextValue "D35" "7.215"
$ws.Range("E35").Value = "  -3.93%  "
Set-TextValue "D36" "11.67"
$ws.Range("E36").Value = "  +4.77%  "
Set-TextValue "D37" "0.2751"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("E38").Value = "  -1.20%  "
Set-TextValue "D39" "14.50"
$ws.Range("E39").Value = "  -0.64%  "
Set-TextValue "D40" "0.09198"
$ws.Range("E40").Value = "  +3.57%  "
Set-TextValue "D41" "0.02730"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("E42").Value = "  +1.15%  "
Set-TextValue "D43" "0.7656"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("E44").Value = "  +4.36%  "
Set-TextValue "D45" "2.606"
$ws.Range("E45").Value = "  +7.09%  "
Set-TextValue "D46" "0.7188"
$ws.Range("E46").Value = "  +1.13%  "
Set-TextValue "D47" "4.220"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("E48").Value = "  +0.50%  "
Set-TextValue "D49" "1.341"
$ws.Range("E49").Value = "  +7.35%  "
Set-TextValue "D50" "140.52"
$ws.Range("E50").Value = "  +0.79%  "
Set-TextValue "D51" "0.07978"
$ws.Range("E51").Value = "  +1.09%  "
